# Add slide transitions to every slide of the deck, matching the
# transition scheme that was applied in the target commit:
#
#   Slide 1 -> peelOff (slow, 1.25s)  [engine has no peelOff/p15:prstTrans
#                                      support, so we fall back to the
#                                      same "fade" that real PowerPoint
#                                      writes as the mc:Fallback branch]
#   Slide 2 -> push, up (slow)
#   Slide 3 -> wipe (slow)
#   Slide 4 -> split, vertical (slow, 1.5s)
#   Slide 5 -> flash (slow)           [engine has no p14:flash support,
#                                      fall back to "fade", matching the
#                                      mc:Fallback branch]
#   Slide 6 -> cover (slow)
#   Slide 7 -> pull (medium)
#   Slide 8 -> circle (slow, 0.8s)
#   Slide 9 -> randomBar, vertical (slow)

$p = $ppt.ActivePresentation

# --- Slide 1: title slide -> peelOff w/ fade fallback (slow, 1.25s) ---
$s1 = $p.Slides.Item(1)
$s1.SlideShowTransition.EntryEffect = 1793   # ppEffectFade
$s1.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 2: "In this presentation" -> push up (slow) ---
$s2 = $p.Slides.Item(2)
$s2.SlideShowTransition.EntryEffect = 3852   # ppEffectPushUp (side push)
$s2.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 3: "Our team" -> wipe (slow) ---
$s3 = $p.Slides.Item(3)
$s3.SlideShowTransition.EntryEffect = 2817   # ppEffectWipeRight (wipe)
$s3.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 4: "Used technologies" -> split vertical (slow, 1.5s) ---
$s4 = $p.Slides.Item(4)
$s4.SlideShowTransition.EntryEffect = 3585   # ppEffectSplitVerticalIn (split)
$s4.SlideShowTransition.Duration = 1.5
$s4.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 5: "Our game" -> flash w/ fade fallback (slow) ---
$s5 = $p.Slides.Item(5)
$s5.SlideShowTransition.EntryEffect = 1793   # ppEffectFade
$s5.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 6: "Stages of realization" -> cover (slow) ---
$s6 = $p.Slides.Item(6)
$s6.SlideShowTransition.EntryEffect = 1281   # ppEffectCoverDown (cover)
$s6.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 7: "Let us show You our repository!" -> pull (medium) ---
$s7 = $p.Slides.Item(7)
$s7.SlideShowTransition.EntryEffect = 2049   # ppEffectPullDown (pull)
$s7.SlideShowTransition.Speed = 2            # ppTransitionSpeedMedium

# --- Slide 8: "Now look at our game!" -> circle (slow, 0.8s) ---
$s8 = $p.Slides.Item(8)
$s8.SlideShowTransition.EntryEffect = 3845   # ppEffectCircle
$s8.SlideShowTransition.Duration = 0.8
$s8.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow

# --- Slide 9: "Thanks for Your attention!" -> randomBar vertical (slow) ---
$s9 = $p.Slides.Item(9)
$s9.SlideShowTransition.EntryEffect = 2305   # ppEffectRandomBarsVertical
$s9.SlideShowTransition.Speed = 1            # ppTransitionSpeedSlow
